$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append two new daily rows (2025-10-25) for both charging stations,
# following the same layout/formatting as the existing data (row 49).
$ws.Range("A49:F49").Copy()
$ws.Range("A50:F50").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A51:F51").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A50").Value = 45955
$ws.Range("B50").Value = "四方坪站"
$ws.Range("C50").Value = 10409.4
$ws.Range("D50").Value = 8571.73
$ws.Range("E50").Value = 3654.9
$ws.Range("F50").Value = 425

$ws.Range("A51").Value = 45955
$ws.Range("B51").Value = "高岭站"
$ws.Range("C51").Value = 5114.2700000000004
$ws.Range("D51").Value = 4154.99
$ws.Range("E51").Value = 1334.4
$ws.Range("F51").Value = 163

# Update the active selection to mirror the edited workbook's last state
$ws.Range("I55").Select()
